$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.902.96"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("D3").Value = "3.278.40"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'580.76"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").Value = "'182.98"
$ws.Range("E6").Value = "  +6.61%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.604"
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").Value = "'0.134"
$ws.Range("E9").Value = "  +7.38%  "
$ws.Range("D10").Value = "'6.72"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("E11").Value = "  +5.87%  "
$ws.Range("D12").Value = "3.849.48"
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("E14").Value = "  +5.35%  "
$ws.Range("D15").Value = "67.865.85"
$ws.Range("E15").Value = "  +3.34%  "
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("D17").Value = "3.284.71"
$ws.Range("E17").Value = "  +3.39%  "
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("E19").Value = "  +4.69%  "
$ws.Range("D20").Value = "'377.57"
$ws.Range("E20").Value = "  +4.43%  "
$ws.Range("E21").Value = "  +5.38%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'71.35"
$ws.Range("E23").Value = "  +3.67%  "
$ws.Range("E24").Value = "  +3.64%  "
$ws.Range("E25").Value = "  +4.91%  "
$ws.Range("D26").Value = "'9.73"
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("E27").Value = "  +3.06%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  +2.71%  "
$ws.Range("D30").Value = "'5.73"
$ws.Range("E30").Value = "  +6.01%  "
$ws.Range("E31").Value = "  +3.92%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.28"
$ws.Range("E32").Value = "  +7.04%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "'0.998"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "'6.98"
$ws.Range("E34").Value = "  +5.36%  "
$ws.Range("E35").Value = "  +5.64%  "
$ws.Range("D36").Value = "'161.83"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").Value = "'0.855"
$ws.Range("E38").Value = "  +2.83%  "
$ws.Range("D39").Value = "'27.03"
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("E40").Value = "  +10.80%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'2.67"
$ws.Range("E41").Value = "  +6.43%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'4.64"
$ws.Range("E42").Value = "  +10.39%  "
$ws.Range("D43").Value = "'25.89"
$ws.Range("E43").Value = "  +8.46%  "
$ws.Range("D44").Value = "'351.54"
$ws.Range("D45").Value = "2.655.62"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("E46").Value = "  +2.87%  "
$ws.Range("D48").Value = "'0.0285"
$ws.Range("E48").Value = "  +4.21%  "
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("E50").Value = "  +5.31%  "
$ws.Range("D51").Value = "'31.11"
$ws.Range("E51").Value = "  +3.80%  "
